$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 65.38041981103575
$ws.Range("C2").Value = 35.858458423859474
$ws.Range("D2").Value = 64.898004515292556
$ws.Range("E2").Value = 39.514648813360758

$ws.Range("B3").Value = 61.871043494024434
$ws.Range("C3").Value = 42.42826502455631
$ws.Range("D3").Value = 53.288661504018926
$ws.Range("E3").Value = 49.837777848804549

$ws.Range("B1:E3").Select()
